$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1428734.2
$ws.Range("I6").Value = 2197873.8
$ws.Range("J6").Value = 332.2857
$ws.Range("K6").Value = 6593621.399999999
$ws.Range("L6").Value = 996.8571000000001
$ws.Range("M6").Value = -6593509.399999999
$ws.Range("N6").Value = -1220.8571
$ws.Range("H100").Value = 2798
$ws.Range("I100").Value = 2120.8
$ws.Range("J100").Value = 3926.6667
$ws.Range("K100").Value = 2120.8
$ws.Range("L100").Value = 3926.6667
$ws.Range("M100").Value = -1579.8
$ws.Range("N100").Value = -5008.6667
$ws.Range("H125").Value = 694809.25
$ws.Range("I125").Value = 276.66666
$ws.Range("J125").Value = 855086
$ws.Range("K125").Value = 2489.99994
$ws.Range("L125").Value = 7695774
$ws.Range("M125").Value = -29.9999399999997
$ws.Range("N125").Value = -7700694
$ws.Range("H131").Value = 3058.1035
$ws.Range("I131").Value = 707.0833
$ws.Range("J131").Value = 4717.647
$ws.Range("K131").Value = 2121.2499
$ws.Range("L131").Value = 14152.941
$ws.Range("M131").Value = 2918.7501
$ws.Range("N131").Value = -24232.941
$ws.Range("H137").Value = 2734.9
$ws.Range("J137").Value = 1410
$ws.Range("L137").Value = 4230
$ws.Range("N137").Value = -9330
$ws.Range("H138").Value = 3267.0908
$ws.Range("I138").Value = 2185
$ws.Range("K138").Value = 6555
$ws.Range("M138").Value = -1415

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13708.53
$ws.Range("I32").Value = 10546.46
$ws.Range("J32").Value = 34870.08
$ws.Range("K32").Value = 10546.46
$ws.Range("L32").Value = 34870.08
$ws.Range("M32").Value = -10259.46
$ws.Range("N32").Value = -35444.08
$ws.Range("H45").Value = 856.3043
$ws.Range("I45").Value = 833.2222
$ws.Range("J45").Value = 939.4
$ws.Range("K45").Value = 833.2222
$ws.Range("L45").Value = 939.4
$ws.Range("M45").Value = -456.2222
$ws.Range("N45").Value = -1693.4
$ws.Range("H92").Value = 31950
$ws.Range("J92").Value = 31950
$ws.Range("L92").Value = 31950
$ws.Range("N92").Value = -36942
$ws.Range("H97").Value = 576.2727
$ws.Range("I97").Value = 448.77777
$ws.Range("J97").Value = 1150
$ws.Range("K97").Value = 448.77777
$ws.Range("L97").Value = 1150
$ws.Range("M97").Value = 47.22223000000002
$ws.Range("N97").Value = -2142
$ws.Range("H102").Value = 2127.7273
$ws.Range("I102").Value = 1515.625
$ws.Range("J102").Value = 3760
$ws.Range("K102").Value = 1515.625
$ws.Range("L102").Value = 3760
$ws.Range("M102").Value = 106.375
$ws.Range("N102").Value = -7004
$ws.Range("H132").Value = 5950.8115
$ws.Range("I132").Value = 4928.706
$ws.Range("J132").Value = 7779.8423
$ws.Range("K132").Value = 14786.118
$ws.Range("L132").Value = 23339.5269
$ws.Range("M132").Value = -12256.118
$ws.Range("N132").Value = -28399.5269

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1425.8948
$ws.Range("I86").Value = 1360.3846
$ws.Range("J86").Value = 1567.8334
$ws.Range("K86").Value = 1360.3846
$ws.Range("L86").Value = 1567.8334
$ws.Range("M86").Value = -237.3846000000001
$ws.Range("N86").Value = -3813.8334
$ws.Range("H89").Value = 1425.8948
$ws.Range("I89").Value = 1360.3846
$ws.Range("J89").Value = 1567.8334
$ws.Range("K89").Value = 6801.923000000001
$ws.Range("L89").Value = 7839.166999999999
$ws.Range("M89").Value = -1185.923000000001
$ws.Range("N89").Value = -19071.167
$ws.Range("H92").Value = 39000
$ws.Range("J92").Value = 39000
$ws.Range("L92").Value = 39000
$ws.Range("N92").Value = -43992
$ws.Range("H94").Value = 1325.0938
$ws.Range("I94").Value = 1305.125
$ws.Range("J94").Value = 1385
$ws.Range("K94").Value = 1305.125
$ws.Range("L94").Value = 1385
$ws.Range("M94").Value = -854.125
$ws.Range("N94").Value = -2287
$ws.Range("H99").Value = 1486.9722
$ws.Range("I99").Value = 1082.9412
$ws.Range("J99").Value = 1848.4736
$ws.Range("K99").Value = 1082.9412
$ws.Range("L99").Value = 1848.4736
$ws.Range("M99").Value = 415.0588
$ws.Range("N99").Value = -4844.4736
$ws.Range("H107").Value = 3356.1538
$ws.Range("I107").Value = 3831.4285
$ws.Range("J107").Value = 1360
$ws.Range("K107").Value = 3831.4285
$ws.Range("L107").Value = 1360
$ws.Range("M107").Value = -1911.4285
$ws.Range("N107").Value = -5200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 40004820
$ws.Range("I86").Value = 52635680
$ws.Range("J86").Value = 7100
$ws.Range("K86").Value = 52635680
$ws.Range("L86").Value = 7100
$ws.Range("M86").Value = -52634557
$ws.Range("N86").Value = -9346
$ws.Range("H89").Value = 40004820
$ws.Range("I89").Value = 52635680
$ws.Range("J89").Value = 7100
$ws.Range("K89").Value = 263178400
$ws.Range("L89").Value = 35500
$ws.Range("M89").Value = -263172784
$ws.Range("N89").Value = -46732
$ws.Range("H92").Value = 30820.2
$ws.Range("J92").Value = 30820.2
$ws.Range("L92").Value = 30820.2
$ws.Range("N92").Value = -35812.2
$ws.Range("H99").Value = 2658.0417
$ws.Range("I99").Value = 2198
$ws.Range("J99").Value = 2986.6428
$ws.Range("K99").Value = 2198
$ws.Range("L99").Value = 2986.6428
$ws.Range("M99").Value = -700
$ws.Range("N99").Value = -5982.6428
$ws.Range("H122").Value = 52632756
$ws.Range("I122").Value = 66667532
$ws.Range("J122").Value = 2351
$ws.Range("K122").Value = 200002596
$ws.Range("L122").Value = 7053
$ws.Range("M122").Value = -200000146
$ws.Range("N122").Value = -11953
$ws.Range("H126").Value = 2658.0417
$ws.Range("I126").Value = 2198
$ws.Range("J126").Value = 2986.6428
$ws.Range("K126").Value = 6594
$ws.Range("L126").Value = 8959.928400000001
$ws.Range("M126").Value = -4124
$ws.Range("N126").Value = -13899.9284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 726.55817
$ws.Range("I5").Value = 499.07693
$ws.Range("J5").Value = 1074.4706
$ws.Range("K5").Value = 1497.23079
$ws.Range("L5").Value = 3223.4118
$ws.Range("M5").Value = -1385.23079
$ws.Range("N5").Value = -3447.4118
$ws.Range("H135").Value = 726.55817
$ws.Range("I135").Value = 499.07693
$ws.Range("J135").Value = 1074.4706
$ws.Range("K135").Value = 4491.69237
$ws.Range("L135").Value = 9670.235400000001
$ws.Range("M135").Value = -1956.69237
$ws.Range("N135").Value = -14740.2354

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 57.363636
$ws.Range("I2").Value = 51.666668
$ws.Range("J2").Value = 64.2
$ws.Range("K2").Value = 51.666668
$ws.Range("L2").Value = 64.2
$ws.Range("M2").Value = 61.333332
$ws.Range("N2").Value = -290.2
$ws.Range("H70").Value = 6206.8887
$ws.Range("I70").Value = 7342.857
$ws.Range("J70").Value = 5484
$ws.Range("K70").Value = 7342.857
$ws.Range("L70").Value = 5484
$ws.Range("M70").Value = -7072.857
$ws.Range("N70").Value = -6024
$ws.Range("H73").Value = 6206.8887
$ws.Range("I73").Value = 7342.857
$ws.Range("J73").Value = 5484
$ws.Range("K73").Value = 7342.857
$ws.Range("L73").Value = 5484
$ws.Range("M73").Value = -6406.857
$ws.Range("N73").Value = -7356
$ws.Range("H80").Value = 2909.4546
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 2900.4
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 2900.4
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -4896.4
$ws.Range("H83").Value = 2909.4546
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 2900.4
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 14502
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -24486
$ws.Range("H92").Value = 6003.5835
$ws.Range("J92").Value = 6003.5835
$ws.Range("L92").Value = 6003.5835
$ws.Range("N92").Value = -9747.583500000001
$ws.Range("H122").Value = 4154.2666
$ws.Range("I122").Value = 5285.6665
$ws.Range("J122").Value = 3400
$ws.Range("K122").Value = 15856.9995
$ws.Range("L122").Value = 10200
$ws.Range("M122").Value = -13406.9995
$ws.Range("N122").Value = -15100
$ws.Range("H132").Value = 4252.5264
$ws.Range("I132").Value = 10012
$ws.Range("J132").Value = 3932.5557
$ws.Range("K132").Value = 30036
$ws.Range("L132").Value = 11797.6671
$ws.Range("M132").Value = -27506
$ws.Range("N132").Value = -16857.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 3333834.8
$ws.Range("I10").Value = 5000250
$ws.Range("K10").Value = 5000250
$ws.Range("M10").Value = -5000110
$ws.Range("H68").Value = 2695.125
$ws.Range("I68").Value = 2266.8333
$ws.Range("J68").Value = 3980
$ws.Range("K68").Value = 2266.8333
$ws.Range("L68").Value = 3980
$ws.Range("M68").Value = -1517.8333
$ws.Range("N68").Value = -5478
$ws.Range("H71").Value = 2695.125
$ws.Range("I71").Value = 2266.8333
$ws.Range("J71").Value = 3980
$ws.Range("K71").Value = 11334.1665
$ws.Range("L71").Value = 19900
$ws.Range("M71").Value = -7590.166499999999
$ws.Range("N71").Value = -27388
$ws.Range("H82").Value = 1775.875
$ws.Range("I82").Value = 1312.3334
$ws.Range("J82").Value = 2371.8572
$ws.Range("K82").Value = 1312.3334
$ws.Range("L82").Value = 2371.8572
$ws.Range("M82").Value = -951.3334
$ws.Range("N82").Value = -3093.8572
$ws.Range("H85").Value = 1775.875
$ws.Range("I85").Value = 1312.3334
$ws.Range("J85").Value = 2371.8572
$ws.Range("K85").Value = 1312.3334
$ws.Range("L85").Value = 2371.8572
$ws.Range("M85").Value = -64.33339999999998
$ws.Range("N85").Value = -4867.8572
$ws.Range("H93").Value = 2400.8572
$ws.Range("I93").Value = 1535.3334
$ws.Range("J93").Value = 3050
$ws.Range("K93").Value = 1535.3334
$ws.Range("L93").Value = 3050
$ws.Range("M93").Value = -287.3334
$ws.Range("N93").Value = -5546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2144.875
$ws.Range("I122").Value = 2384.8333
$ws.Range("J122").Value = 1425
$ws.Range("K122").Value = 7154.499899999999
$ws.Range("L122").Value = 4275
$ws.Range("M122").Value = -4704.499899999999
$ws.Range("N122").Value = -9175
$ws.Range("H132").Value = 1628.3334
$ws.Range("I132").Value = 1358.9412
$ws.Range("J132").Value = 1980.6154
$ws.Range("K132").Value = 4076.8236
$ws.Range("L132").Value = 5941.8462
$ws.Range("M132").Value = -1546.8236
$ws.Range("N132").Value = -11001.8462

Write-Output "Updated 273 cells."